# The edit re-orders the data rows (2-16) of the "Artfynd" sheet: each
# row's full content (columns A:AY) is relocated to a different row
# position (a pure permutation of whole rows - no cell values are
# otherwise changed).
#
# Mapping of original row number -> new row number (derived from the Id
# column in A, which travels together with the rest of that row's data):
#   2->7, 3->2, 4->8, 5->9, 6->3, 7->10, 8->11, 9->4, 10->12,
#   11->13, 12->14, 13->15, 14->16, 15->5, 16->6

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstRow = 2
$lastRow = 16

$mapping = @{
    2  = 7
    3  = 2
    4  = 8
    5  = 9
    6  = 3
    7  = 10
    8  = 11
    9  = 4
    10 = 12
    11 = 13
    12 = 14
    13 = 15
    14 = 16
    15 = 5
    16 = 6
}

# Every row already carries an (empty) inline-string placeholder in
# columns I, AT and AY. Those columns are left out of the bulk copy
# below (writing "" through COM would delete the cell outright instead
# of keeping the empty-string marker), so the segments here purposely
# skip I (col 9) and AT (col 46); AY (col 51) is handled separately
# afterwards since its text content (when present) does need to move
# with the row.
$segments = @(
    @{ Start = "A";  End = "H"  },
    @{ Start = "J";  End = "AS" },
    @{ Start = "AU"; End = "AX" }
)

# Columns whose text looks like a date (YYYY-MM-DD) - written back with
# a leading apostrophe plus a forced General format so Excel keeps
# storing the literal text instead of re-typing it as a real date.
$dateCols = @("Y", "AA")

# 1) Snapshot every source row (values only) before any writes happen.
$rowsData = @{}
$ayData = @{}
for ($r = $firstRow; $r -le $lastRow; $r++) {
    $rowsData[$r] = @{}
    foreach ($seg in $segments) {
        $rng = $ws.Range("$($seg.Start)$r`:$($seg.End)$r")
        $rowsData[$r][$seg.Start] = $rng.Value2
    }
    $ayData[$r] = $ws.Range("AY$r").Value2
}

# 2) Write each snapshot into its destination row.
foreach ($srcRow in $rowsData.Keys) {
    $destRow = $mapping[$srcRow]

    foreach ($seg in $segments) {
        $destRng = $ws.Range("$($seg.Start)$destRow`:$($seg.End)$destRow")
        $destRng.Value2 = $rowsData[$srcRow][$seg.Start]
    }

    # AY: only touch the destination cell when the source actually had
    # text, so rows that came from an empty AY keep the destination's
    # own pre-existing empty marker untouched.
    $ayVal = $ayData[$srcRow]
    if ($null -ne $ayVal -and $ayVal -ne "") {
        $ws.Range("AY$destRow").Value2 = $ayVal
    }

    # Re-apply the date-like text columns explicitly with a leading
    # apostrophe (plus General format) so Excel keeps them as plain
    # text instead of silently re-typing them as dates.
    foreach ($col in $dateCols) {
        $text = $ws.Range("$col$destRow").Value2
        if ($null -ne $text -and $text -ne "") {
            $ws.Range("$col$destRow").Value = "'" + $text
            $ws.Range("$col$destRow").NumberFormat = "General"
        }
    }
}
